$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$rows = @(
    @("2026-02-01", "11:40:23", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:40:23", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:40:23", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:40:23", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:40:23", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:40:29", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:40:40", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:40:50", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:41:01", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:41:11", "11:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "11:41:21", "11:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 80
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Column A holds a date-like string ("2026-02-01"). Excel's COM layer
    # auto-converts such strings to real dates, so force text formatting,
    # assign the value, then clear the formatting again so the cell keeps
    # no explicit style (matching the rest of the sheet) while its content
    # stays the literal text string.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[0]
    $cellA.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}
